# "Water Tower Design Base" - update FIRE and WATER tower stat sheets.

$wb = $excel.ActiveWorkbook
$wsFire  = $wb.Worksheets.Item("FIRE")
$wsWater = $wb.Worksheets.Item("WATER")

# ---------------------------------------------------------------------------
# FIRE sheet: reduce ability-stat (column F) for rows 21-23 from 1.2 to 0.8
# ---------------------------------------------------------------------------
$wsFire.Range("F21").Value = 0.8
$wsFire.Range("F22").Value = 0.8
$wsFire.Range("F23").Value = 0.8

# ---------------------------------------------------------------------------
# WATER sheet: update attacktime (H), type/count (I), projspf (J) and
# projspeed (K) columns, plus a couple of ability (F) tweaks.
# Values are assigned top-to-bottom so new shared strings are appended in
# the same order as the target workbook.
# ---------------------------------------------------------------------------

# F column ability tweaks
$wsWater.Range("F2").Value  = 0.6
$wsWater.Range("F3").Value  = 0.6
$wsWater.Range("F22").Value = 0.5
$wsWater.Range("F23").Value = 0.5

# H column (attacktime) new text values
$wsWater.Range("H2").Value  = "0.2;0.8;"
$wsWater.Range("H3").Value  = "0.2;0.9;"
$wsWater.Range("H4").Value  = "0.3;"
$wsWater.Range("H5").Value  = "0.1;"
$wsWater.Range("H6").Value  = "0.5;"
$wsWater.Range("H7").Value  = "0.4;"
$wsWater.Range("H8").Value  = "0.4;"
$wsWater.Range("H9").Value  = "0.4;"
$wsWater.Range("H10").Value = "0.6;"
$wsWater.Range("H11").Value = "0.6;"
$wsWater.Range("H12").Value = "0.6;"
$wsWater.Range("H13").Value = "0.3;"
$wsWater.Range("H14").Value = "0.3;"
$wsWater.Range("H15").Value = "0.5;"
$wsWater.Range("H16").Value = "0.5;"
$wsWater.Range("H17").Value = "0.5;"
$wsWater.Range("H18").Value = "0.6;"
$wsWater.Range("H19").Value = "0.6;"
$wsWater.Range("H20").Value = "0.6;"
$wsWater.Range("H21").Value = "0.6;"
$wsWater.Range("H22").Value = "0.2;0.8;1.4;"
$wsWater.Range("H23").Value = "0.2;0.8;1.4;"

# I/J/K columns: projectile type/count + spf/speed for the rows that gained
# projectile data in this revision
$wsWater.Range("I5").Value  = 1
$wsWater.Range("J5").Value  = 0.1
$wsWater.Range("K5").Value  = 0.3

$wsWater.Range("I6").Value  = 1
$wsWater.Range("J6").Value  = 0.1
$wsWater.Range("K6").Value  = 0.3

$wsWater.Range("I10").Value = 1
$wsWater.Range("J10").Value = 0.1
$wsWater.Range("K10").Value = 0.3

$wsWater.Range("I11").Value = 1
$wsWater.Range("J11").Value = 0.1
$wsWater.Range("K11").Value = 0.3

$wsWater.Range("I12").Value = 1
$wsWater.Range("J12").Value = 0.1
$wsWater.Range("K12").Value = 0.3

$wsWater.Range("I15").Value = 1
$wsWater.Range("J15").Value = 0.1
$wsWater.Range("K15").Value = 0.3

$wsWater.Range("I16").Value = 1
$wsWater.Range("J16").Value = 0.1
$wsWater.Range("K16").Value = 0.3

$wsWater.Range("I17").Value = 1
$wsWater.Range("J17").Value = 0.1
$wsWater.Range("K17").Value = 0.3

$wsWater.Range("I18").Value = 2
$wsWater.Range("J18").Value = 0.05
$wsWater.Range("K18").Value = 0.45

$wsWater.Range("I19").Value = 2
$wsWater.Range("J19").Value = 0.05
$wsWater.Range("K19").Value = 0.45

$wsWater.Range("I20").Value = 2
$wsWater.Range("J20").Value = 0.05
$wsWater.Range("K20").Value = 0.45

$wsWater.Range("I21").Value = 2
$wsWater.Range("J21").Value = 0.05
$wsWater.Range("K21").Value = 0.45

$wsWater.Range("I22").Value = 1
$wsWater.Range("J22").Value = 0.1
$wsWater.Range("K22").Value = 0.3

$wsWater.Range("I23").Value = 1
$wsWater.Range("J23").Value = 0.1
$wsWater.Range("K23").Value = 0.3

# ---------------------------------------------------------------------------
# View state: WATER becomes the active / selected tab, with its own
# selection, and FIRE keeps a plain (non-active) selection.
# ---------------------------------------------------------------------------
$wsFire.Activate()
$wsFire.Range("F24").Select()

$wsWater.Activate()
$wsWater.Range("K24").Select()

$wb.Save()
